$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column N (14) - this shifts the existing
# "Description" column (and its data) from N to O.
$ws.Columns.Item(14).Insert()

# New column header
$ws.Cells.Item(1, 14).Value = "eIDAS RequesterID"

# New column data values - set in the same order the original authors did
# (row 6 first, then rows 2-5, then rows 7-8) so shared-string allocation
# order matches.
$ws.Cells.Item(6, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298e"
$ws.Cells.Item(2, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298a"
$ws.Cells.Item(3, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298b"
$ws.Cells.Item(4, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298c"
$ws.Cells.Item(5, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298d"
$ws.Cells.Item(7, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298f"
$ws.Cells.Item(8, 14).Value = "58ee2267-7864-4e09-958b-b53c3135298g"

# Leave the selection where the author ended up after typing the last value
# (one row below the data, still in the new column).
$ws.Cells.Item(9, 14).Select() | Out-Null
